# Update the date line and the twenty-five three-digit ÷ one-digit
# division answers in the practice table to the new set of problems.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-03-23 Sunday" "2025-03-24 Monday"

Replace-Text "541÷7=77, 2"   "639÷2=319, 1"
Replace-Text "694÷4=173, 2"  "108÷2=54, 0"
Replace-Text "680÷3=226, 2"  "378÷9=42, 0"
Replace-Text "340÷7=48, 4"   "519÷3=173, 0"
Replace-Text "228÷3=76, 0"   "818÷8=102, 2"

Replace-Text "741÷9=82, 3"   "206÷4=51, 2"
Replace-Text "825÷4=206, 1"  "235÷2=117, 1"
Replace-Text "738÷3=246, 0"  "271÷5=54, 1"
Replace-Text "463÷4=115, 3"  "330÷4=82, 2"
Replace-Text "877÷2=438, 1"  "798÷7=114, 0"

Replace-Text "571÷3=190, 1"  "556÷4=139, 0"
Replace-Text "560÷4=140, 0"  "572÷8=71, 4"
Replace-Text "711÷2=355, 1"  "200÷3=66, 2"
Replace-Text "193÷5=38, 3"   "981÷7=140, 1"
Replace-Text "663÷8=82, 7"   "315÷8=39, 3"

Replace-Text "751÷4=187, 3"  "808÷6=134, 4"
Replace-Text "498÷7=71, 1"   "980÷6=163, 2"
Replace-Text "894÷7=127, 5"  "857÷8=107, 1"
Replace-Text "432÷2=216, 0"  "551÷6=91, 5"
Replace-Text "310÷4=77, 2"   "489÷8=61, 1"

Replace-Text "668÷5=133, 3"  "902÷7=128, 6"
Replace-Text "428÷8=53, 4"   "534÷4=133, 2"
Replace-Text "717÷6=119, 3"  "186÷2=93, 0"
Replace-Text "653÷4=163, 1"  "393÷6=65, 3"
Replace-Text "287÷9=31, 8"   "560÷7=80, 0"

Write-Output "done"
